# Append: 2025-11-23 18:31 JST
# Update the "取得日時" (fetched-at) timestamp in column A for rows 2-6
# on the active sheet ("ランサーズ") from 2025-11-23 18:23:41 to 2025-11-23 18:31:21.
# Values are stored as plain text, so we set them explicitly as strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-23 18:31:21"

foreach ($r in 2..6) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
